$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1465.55
$ws.Range("I19").Value = 2618.2
$ws.Range("J19").Value = 1081.3334
$ws.Range("K19").Value = 2618.2
$ws.Range("L19").Value = 1081.3334
$ws.Range("M19").Value = -2443.2
$ws.Range("N19").Value = -1431.3334
$ws.Range("H28").Value = 2667.7856
$ws.Range("I28").Value = 2568.75
$ws.Range("K28").Value = 2568.75
$ws.Range("M28").Value = -2083.75
$ws.Range("H76").Value = 6561.5
$ws.Range("I76").Value = 5388.2354
$ws.Range("J76").Value = 8777.666999999999
$ws.Range("K76").Value = 5388.2354
$ws.Range("L76").Value = 8777.666999999999
$ws.Range("M76").Value = -5073.2354
$ws.Range("N76").Value = -9407.666999999999
$ws.Range("H79").Value = 6561.5
$ws.Range("I79").Value = 5388.2354
$ws.Range("J79").Value = 8777.666999999999
$ws.Range("K79").Value = 5388.2354
$ws.Range("L79").Value = 8777.666999999999
$ws.Range("M79").Value = -4296.2354
$ws.Range("N79").Value = -10961.667
$ws.Range("H97").Value = 900
$ws.Range("J97").Value = 900
$ws.Range("L97").Value = 2700
$ws.Range("N97").Value = -3692
$ws.Range("H106").Value = 629350
$ws.Range("I106").Value = 718685.7
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 718685.7
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -718054.7
$ws.Range("N106").Value = -5262
$ws.Range("H138").Value = 2111.8364
$ws.Range("I138").Value = 1250.0625
$ws.Range("J138").Value = 3310.8262
$ws.Range("K138").Value = 3750.1875
$ws.Range("L138").Value = 9932.4786
$ws.Range("M138").Value = 1389.8125
$ws.Range("N138").Value = -20212.4786
$ws.Range("H141").Value = 1588.871
$ws.Range("I141").Value = 1087.2222
$ws.Range("J141").Value = 4975
$ws.Range("K141").Value = 3261.6666
$ws.Range("L141").Value = 14925
$ws.Range("M141").Value = 1918.3334
$ws.Range("N141").Value = -25285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1303.2693
$ws.Range("I45").Value = 1350
$ws.Range("J45").Value = 1147.5
$ws.Range("K45").Value = 1350
$ws.Range("L45").Value = 1147.5
$ws.Range("M45").Value = -973
$ws.Range("N45").Value = -1901.5
$ws.Range("H88").Value = 2990.9092
$ws.Range("I88").Value = 2871.4285
$ws.Range("J88").Value = 3200
$ws.Range("K88").Value = 2871.4285
$ws.Range("L88").Value = 3200
$ws.Range("M88").Value = -2465.4285
$ws.Range("N88").Value = -4012
$ws.Range("H91").Value = 2990.9092
$ws.Range("I91").Value = 2871.4285
$ws.Range("J91").Value = 3200
$ws.Range("K91").Value = 2871.4285
$ws.Range("L91").Value = 3200
$ws.Range("M91").Value = -1467.4285
$ws.Range("N91").Value = -6008
$ws.Range("H97").Value = 5020.222
$ws.Range("I97").Value = 5648.8423
$ws.Range("K97").Value = 5648.8423
$ws.Range("M97").Value = -5152.8423
$ws.Range("H102").Value = 2308.5715
$ws.Range("I102").Value = 2505
$ws.Range("J102").Value = 2230
$ws.Range("K102").Value = 2505
$ws.Range("L102").Value = 2230
$ws.Range("M102").Value = -883
$ws.Range("N102").Value = -5474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 71431120
$ws.Range("I86").Value = 2636.3635
$ws.Range("J86").Value = 150002450
$ws.Range("K86").Value = 2636.3635
$ws.Range("L86").Value = 150002450
$ws.Range("M86").Value = -1513.3635
$ws.Range("N86").Value = -150004696
$ws.Range("H89").Value = 71431120
$ws.Range("I89").Value = 2636.3635
$ws.Range("J89").Value = 150002450
$ws.Range("K89").Value = 13181.8175
$ws.Range("L89").Value = 750012250
$ws.Range("M89").Value = -7565.817499999999
$ws.Range("N89").Value = -750023482
$ws.Range("H94").Value = 1092
$ws.Range("I94").Value = 999.7692
$ws.Range("J94").Value = 1211.9
$ws.Range("K94").Value = 999.7692
$ws.Range("L94").Value = 1211.9
$ws.Range("M94").Value = -548.7692
$ws.Range("N94").Value = -2113.9
$ws.Range("H99").Value = 782.6667
$ws.Range("I99").Value = 759
$ws.Range("J99").Value = 830
$ws.Range("K99").Value = 759
$ws.Range("L99").Value = 830
$ws.Range("M99").Value = 739
$ws.Range("N99").Value = -3826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 38745.5
$ws.Range("J140").Value = 38745.5
$ws.Range("L140").Value = 38745.5
$ws.Range("N140").Value = -49105.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4841.25
$ws.Range("I133").Value = 3639
$ws.Range("J133").Value = 5700
$ws.Range("K133").Value = 10917
$ws.Range("L133").Value = 17100
$ws.Range("M133").Value = -5857
$ws.Range("N133").Value = -27220
$ws.Range("H137").Value = 4287.1177
$ws.Range("I137").Value = 3423
$ws.Range("J137").Value = 4647.1665
$ws.Range("K137").Value = 10269
$ws.Range("L137").Value = 13941.4995
$ws.Range("M137").Value = -5169
$ws.Range("N137").Value = -24141.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2039
$ws.Range("I97").Value = 2022.7858
$ws.Range("J97").Value = 2084.4
$ws.Range("K97").Value = 2022.7858
$ws.Range("L97").Value = 2084.4
$ws.Range("M97").Value = -1526.7858
$ws.Range("N97").Value = -3076.4
$ws.Range("H126").Value = 3375.0205
$ws.Range("I126").Value = 2487.5356
$ws.Range("J126").Value = 4558.3335
$ws.Range("K126").Value = 7462.6068
$ws.Range("L126").Value = 13675.0005
$ws.Range("M126").Value = -4992.6068
$ws.Range("N126").Value = -18615.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 658.3
$ws.Range("I55").Value = 396.83334
$ws.Range("J55").Value = 1050.5
$ws.Range("K55").Value = 396.83334
$ws.Range("L55").Value = 1050.5
$ws.Range("M55").Value = -223.83334
$ws.Range("N55").Value = -1396.5
$ws.Range("H93").Value = 1586.8572
$ws.Range("I93").Value = 1600
$ws.Range("J93").Value = 1581.6
$ws.Range("K93").Value = 1600
$ws.Range("L93").Value = 1581.6
$ws.Range("M93").Value = -352
$ws.Range("N93").Value = -4077.6
$ws.Range("H122").Value = 6009.6665
$ws.Range("I122").Value = 6382.1333
$ws.Range("J122").Value = 5388.8887
$ws.Range("K122").Value = 19146.3999
$ws.Range("L122").Value = 16166.6661
$ws.Range("M122").Value = -16696.3999
$ws.Range("N122").Value = -21066.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
